$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.644.72"
$ws.Range("E2").Value = "  -2.60%  "
$ws.Range("D3").Value = "3.545.85"
$ws.Range("E3").Value = "  -3.46%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'613.13"
$ws.Range("E5").Value = "  -5.90%  "
$ws.Range("D6").Value = "'153.74"
$ws.Range("E6").Value = "  -3.43%  "
$ws.Range("D7").Value = "3.546.10"
$ws.Range("E7").Value = "  -3.26%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -2.06%  "
$ws.Range("E10").Value = "  -2.33%  "
$ws.Range("D11").Value = "'6.88"
$ws.Range("E11").Value = "  -2.28%  "
$ws.Range("D12").Value = "'0.431"
$ws.Range("E12").Value = "  -1.69%  "
$ws.Range("E13").Value = "  -2.88%  "
$ws.Range("D14").Value = "4.150.30"
$ws.Range("E14").Value = "  -3.34%  "
$ws.Range("D15").Value = "'32.17"
$ws.Range("E15").Value = "  -0.83%  "
$ws.Range("D16").Value = "3.550.23"
$ws.Range("E16").Value = "  -3.48%  "
$ws.Range("D17").Value = "67.724.79"
$ws.Range("E17").Value = "  -2.45%  "
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("D19").Value = "'6.37"
$ws.Range("E19").Value = "  -0.64%  "
$ws.Range("D20").Value = "'15.55"
$ws.Range("E20").Value = "  -2.54%  "
$ws.Range("D21").Value = "'453.38"
$ws.Range("E21").Value = "  -2.59%  "
$ws.Range("D22").Value = "'9.45"
$ws.Range("E22").Value = "  -2.25%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "'77.62"
$ws.Range("E24").Value = "  -2.31%  "
$ws.Range("D25").Value = "3.693.57"
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  -5.77%  "
$ws.Range("D28").Value = "'10.52"
$ws.Range("E28").Value = "  -2.74%  "
$ws.Range("D29").Value = "'8.34"
$ws.Range("E29").Value = "  -6.02%  "
$ws.Range("D31").Value = "'1.62"
$ws.Range("E31").Value = "  -2.54%  "
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("D33").Value = "'25.88"
$ws.Range("E33").Value = "  -2.65%  "
$ws.Range("D34").Value = "'1.90"
$ws.Range("E34").Value = "  -4.13%  "
$ws.Range("D35").Value = "'6.22"
$ws.Range("E35").Value = "  -3.28%  "
$ws.Range("E36").Value = "  -1.95%  "
$ws.Range("D37").Value = "3.550.18"
$ws.Range("E37").Value = "  -3.12%  "
$ws.Range("D38").Value = "'8.05"
$ws.Range("E38").Value = "  -3.83%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("D41").Value = "'177.02"
$ws.Range("E41").Value = "  -1.17%  "
$ws.Range("D42").Value = "'0.0882"
$ws.Range("E42").Value = "  -1.25%  "
$ws.Range("E43").Value = "  -5.57%  "
$ws.Range("D44").Value = "'2.07"
$ws.Range("E44").Value = "  -5.32%  "
$ws.Range("D45").Value = "'0.887"
$ws.Range("E45").Value = "  -4.52%  "
$ws.Range("D46").Value = "'29.05"
$ws.Range("E46").Value = "  +7.35%  "
$ws.Range("D47").Value = "'45.87"
$ws.Range("E47").Value = "  -2.15%  "
$ws.Range("E48").Value = "  -4.42%  "
$ws.Range("E49").Value = "  -1.35%  "
$ws.Range("E50").Value = "  -4.90%  "
$ws.Range("E51").Value = "  -4.13%  "
